$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string shown in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 15:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1795635
$ws.Range("C4").Value = 2105
$ws.Range("D4").Value = 519612
$ws.Range("E4").Value = 1171442
$ws.Range("G4").Value = 39
$ws.Range("H4").Value = 104581

# Row 11 - Alemania
$ws.Range("B11").Value = 183089
$ws.Range("C11").Value = 70
$ws.Range("E11").Value = 9591
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = 8598

# Row 12 - India
$ws.Range("B12").Value = 175434
$ws.Range("C12").Value = 1943
$ws.Range("E12").Value = 87506

# Row 52 - Serbia
$ws.Range("B52").Value = 11381
$ws.Range("C52").Value = 27
$ws.Range("D52").Value = 6606
$ws.Range("E52").Value = 4533

# Row 61 - Noruega
$ws.Range("B61").Value = 8435
$ws.Range("C61").Value = 13
$ws.Range("E61").Value = 472

# Row 76 - Tayikistan
$ws.Range("B76").Value = 3807
$ws.Range("C76").Value = 121
$ws.Range("D76").Value = 1865
$ws.Range("E76").Value = 1895

# Row 94 - Islandia
$ws.Range("B94").Value = 1806
$ws.Range("C94").Value = 1
$ws.Range("E94").Value = 2

# Row 156 - Mozambique
$ws.Range("B156").Value = 244
$ws.Range("C156").Value = 10
$ws.Range("D156").Value = 90
$ws.Range("E156").Value = 152
